# Updated cryptos list — apply value changes per the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    # Force literal text (preserve values like "398.80" or "57.453.75"
    # instead of letting Excel auto-convert look-alike numbers),
    # then clear the temporary number format so the cell keeps its
    # original (default) style, matching the source data.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "57.453.75"
$ws.Range("E2").Value = "  +2.65%  "
Set-TextValue $ws.Range("D3") "3.270.86"
$ws.Range("E3").Value = "  +1.53%  "
$ws.Range("E4").Value = "  -0.07%  "
Set-TextValue $ws.Range("D5") "398.80"
$ws.Range("E5").Value = "  +0.52%  "
Set-TextValue $ws.Range("D6") "109.08"
$ws.Range("E6").Value = "  -1.18%  "
Set-TextValue $ws.Range("D7") "0.580"
$ws.Range("E7").Value = "  +5.09%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +0.50%  "
Set-TextValue $ws.Range("D10") "39.29"
$ws.Range("E10").Value = "  +0.40%  "
Set-TextValue $ws.Range("D11") "0.0962"
$ws.Range("E11").Value = "  +5.81%  "
Set-TextValue $ws.Range("D12") "0.143"
$ws.Range("E12").Value = "  +1.25%  "
Set-TextValue $ws.Range("D13") "3.784.71"
$ws.Range("E13").Value = "  +1.39%  "
Set-TextValue $ws.Range("D14") "8.28"
$ws.Range("E14").Value = "  +2.85%  "
Set-TextValue $ws.Range("D15") "18.97"
$ws.Range("E15").Value = "  -0.11%  "
Set-TextValue $ws.Range("D16") "3.270.61"
$ws.Range("E16").Value = "  +1.23%  "
$ws.Range("E17").Value = "  -1.11%  "
Set-TextValue $ws.Range("D18") "11.15"
$ws.Range("E18").Value = "  +3.24%  "
Set-TextValue $ws.Range("D19") "57.264.17"
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("E21").Value = "  +5.68%  "
$ws.Range("E22").Value = "  -0.34%  "
Set-TextValue $ws.Range("D23") "297.28"
$ws.Range("E23").Value = "  -0.34%  "
Set-TextValue $ws.Range("D24") "74.37"
$ws.Range("E24").Value = "  -1.37%  "
Set-TextValue $ws.Range("D25") "3.19"
$ws.Range("E25").Value = "  -0.70%  "
Set-TextValue $ws.Range("D26") "28.10"
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("B27").Value = "Filecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D27") "7.89"
$ws.Range("E27").Value = "  -3.08%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D28") "4.39"
$ws.Range("E28").Value = "  +0.38%  "
Set-TextValue $ws.Range("D29") "7.42"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("E30").Value = "  -1.41%  "
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("E32").Value = "  +1.47%  "
Set-TextValue $ws.Range("D33") "11.28"
$ws.Range("E33").Value = "  +1.56%  "
Set-TextValue $ws.Range("D34") "40.49"
$ws.Range("E34").Value = "  +12.43%  "
Set-TextValue $ws.Range("D35") "0.0499"
$ws.Range("E35").Value = "  +1.88%  "
$ws.Range("E36").Value = "  +0.98%  "
Set-TextValue $ws.Range("D37") "51.50"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D38") "0.999"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D39") "3.09"
$ws.Range("E39").Value = "  -0.90%  "
Set-TextValue $ws.Range("D40") "3.49"
$ws.Range("E40").Value = "  -0.89%  "
Set-TextValue $ws.Range("D41") "138.04"
$ws.Range("E41").Value = "  +3.28%  "
$ws.Range("E42").Value = "  +2.14%  "
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("E44").Value = "  -1.97%  "
Set-TextValue $ws.Range("D45") "16.87"
$ws.Range("E45").Value = "  -2.52%  "
Set-TextValue $ws.Range("D46") "3.91"
$ws.Range("E46").Value = "  -1.98%  "
Set-TextValue $ws.Range("D47") "22.50"
$ws.Range("E47").Value = "  +1.41%  "
Set-TextValue $ws.Range("D48") "2.22"
$ws.Range("E48").Value = "  +4.70%  "
Set-TextValue $ws.Range("D49") "2.158.54"
$ws.Range("E49").Value = "  +1.42%  "
Set-TextValue $ws.Range("D50") "2.46"
$ws.Range("E50").Value = "  -0.13%  "
Set-TextValue $ws.Range("D51") "1.93"
$ws.Range("E51").Value = "  -11.27%  "
